# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates to the Goblin Profits workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 845.6667
$ws.Range("I15").Value = 845.6667
$ws.Range("K15").Value = 2537.0001
$ws.Range("M15").Value = -2368.0001

# Row 33
$ws.Range("H33").Value = 844.86664
$ws.Range("I33").Value = 208.45454
$ws.Range("J33").Value = 2595
$ws.Range("K33").Value = 208.45454
$ws.Range("L33").Value = 2595
$ws.Range("M33").Value = 20.54545999999999
$ws.Range("N33").Value = -3053

# Row 40
$ws.Range("H40").Value = 2501.4546
$ws.Range("J40").Value = 3222.1428
$ws.Range("L40").Value = 3222.1428
$ws.Range("N40").Value = -3572.1428

# Row 64
$ws.Range("H64").Value = 7942.8286
$ws.Range("J64").Value = 9307.654
$ws.Range("L64").Value = 9307.654
$ws.Range("N64").Value = -9803.654

# Row 67
$ws.Range("H67").Value = 7942.8286
$ws.Range("J67").Value = 9307.654
$ws.Range("L67").Value = 9307.654
$ws.Range("N67").Value = -11023.654

# Row 107
$ws.Range("H107").Value = 1094.7
$ws.Range("I107").Value = 1094.7
$ws.Range("K107").Value = 1094.7
$ws.Range("M107").Value = 825.3

# Row 120
$ws.Range("H120").Value = 114000
$ws.Range("J120").Value = 114000
$ws.Range("L120").Value = 114000
$ws.Range("N120").Value = -123676

# Row 125
$ws.Range("H125").Value = 1774
$ws.Range("J125").Value = 1850
$ws.Range("L125").Value = 16650
$ws.Range("N125").Value = -21570

# Row 137
$ws.Range("H137").Value = 1287.697
$ws.Range("I137").Value = 1134.2307
$ws.Range("K137").Value = 3402.6921
$ws.Range("M137").Value = -852.6921000000002

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4313.9062
$ws.Range("I32").Value = 4255.077
$ws.Range("K32").Value = 4255.077
$ws.Range("M32").Value = -3968.077

# Row 39
$ws.Range("H39").Value = 2000000
$ws.Range("I39").Value = 2000000
$ws.Range("K39").Value = 2000000
$ws.Range("M39").Value = -1999480

# Row 61
$ws.Range("H61").Value = 5722.857
$ws.Range("I61").Value = 5636.25
$ws.Range("K61").Value = 5636.25
$ws.Range("M61").Value = -5424.25

# Row 74
$ws.Range("H74").Value = 2382.1428
$ws.Range("I74").Value = 2029.1666
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 2029.1666
$ws.Range("L74").Value = 4500
$ws.Range("M74").Value = -1155.1666
$ws.Range("N74").Value = -6248

# Row 77
$ws.Range("H77").Value = 2382.1428
$ws.Range("I77").Value = 2029.1666
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 10145.833
$ws.Range("L77").Value = 22500
$ws.Range("M77").Value = -5777.833000000001
$ws.Range("N77").Value = -31236

# Row 102
$ws.Range("H102").Value = 3394.1
$ws.Range("I102").Value = 1742.625
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 1742.625
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -120.625
$ws.Range("N102").Value = -13244

# Row 132
$ws.Range("H132").Value = 1682.9445
$ws.Range("I132").Value = 1713.0588
$ws.Range("K132").Value = 5139.1764
$ws.Range("M132").Value = -2609.1764

# Row 136
$ws.Range("H136").Value = 5722.857
$ws.Range("I136").Value = 5636.25
$ws.Range("K136").Value = 16908.75
$ws.Range("M136").Value = -14358.75

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 3218.4375
$ws.Range("I99").Value = 1700
$ws.Range("J99").Value = 3724.5833
$ws.Range("K99").Value = 1700
$ws.Range("L99").Value = 3724.5833
$ws.Range("M99").Value = -202
$ws.Range("N99").Value = -6720.5833

# Row 105
$ws.Range("H105").Value = 1610.75
$ws.Range("I105").Value = 1484.6538
$ws.Range("K105").Value = 1484.6538
$ws.Range("M105").Value = 262.3462

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3290.5293
$ws.Range("I31").Value = 1441.579
$ws.Range("J31").Value = 5632.533
$ws.Range("K31").Value = 1441.579
$ws.Range("L31").Value = 5632.533
$ws.Range("M31").Value = -1146.579
$ws.Range("N31").Value = -6222.533

# Row 34
$ws.Range("H34").Value = 3290.5293
$ws.Range("I34").Value = 1441.579
$ws.Range("J34").Value = 5632.533
$ws.Range("K34").Value = 1441.579
$ws.Range("L34").Value = 5632.533
$ws.Range("M34").Value = -1239.579
$ws.Range("N34").Value = -6036.533

# Row 48
$ws.Range("H48").Value = 24966.666
$ws.Range("J48").Value = 24966.666
$ws.Range("L48").Value = 24966.666
$ws.Range("N48").Value = -25918.666

# Row 62
$ws.Range("H62").Value = 10299
$ws.Range("I62").Value = 2165
$ws.Range("J62").Value = 22500
$ws.Range("K62").Value = 2165
$ws.Range("L62").Value = 22500
$ws.Range("M62").Value = -1541
$ws.Range("N62").Value = -23748

# Row 65
$ws.Range("H65").Value = 10299
$ws.Range("I65").Value = 2165
$ws.Range("J65").Value = 22500
$ws.Range("K65").Value = 10825
$ws.Range("L65").Value = 112500
$ws.Range("M65").Value = -7705
$ws.Range("N65").Value = -118740

# Row 112
$ws.Range("H112").Value = 75000
$ws.Range("J112").Value = 75000
$ws.Range("L112").Value = 75000
$ws.Range("N112").Value = -77954

# Row 132
$ws.Range("H132").Value = 1630.64
$ws.Range("I132").Value = 1549.6842
$ws.Range("K132").Value = 4649.0526
$ws.Range("M132").Value = -2119.0526

$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 130.38461
$ws.Range("I33").Value = 71.375
$ws.Range("K33").Value = 428.25
$ws.Range("M33").Value = -145.25

# Row 55
$ws.Range("H55").Value = 339533.66
$ws.Range("J55").Value = 6846.5386
$ws.Range("L55").Value = 20539.6158
$ws.Range("N55").Value = -20893.6158

# Row 98
$ws.Range("H98").Value = 124
$ws.Range("J98").Value = 124
$ws.Range("L98").Value = 372
$ws.Range("N98").Value = -3368

$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 1800
$ws.Range("I31").Value = 1800
$ws.Range("K31").Value = 1800
$ws.Range("M31").Value = -1508

# Row 35
$ws.Range("H35").Value = 25000
$ws.Range("I35").Value = 25000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 25000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -24702
$ws.Range("N35").ClearContents()

# Row 37
$ws.Range("H37").Value = 1800
$ws.Range("I37").Value = 1800
$ws.Range("K37").Value = 1800
$ws.Range("M37").Value = -1523

# Row 132
$ws.Range("H132").Value = 2009.8182
$ws.Range("I132").Value = 1711
$ws.Range("K132").Value = 5133
$ws.Range("M132").Value = -2603

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5913.25
$ws.Range("I40").Value = 4473.8887
$ws.Range("K40").Value = 4473.8887
$ws.Range("M40").Value = -4337.8887

# Row 111
$ws.Range("H111").Value = 201750
$ws.Range("J111").Value = 201750
$ws.Range("L111").Value = 201750
$ws.Range("N111").Value = -209930

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 1074
$ws.Range("I5").Value = 1111
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1111
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = -999
$ws.Range("N5").Value = -1224

# Row 46
$ws.Range("H46").Value = 82949.5
$ws.Range("J46").Value = 82949.5
$ws.Range("L46").Value = 82949.5
$ws.Range("N46").Value = -83411.5

# Row 62
$ws.Range("H62").Value = 15799.044
$ws.Range("I62").Value = 5975
$ws.Range("J62").Value = 17867.264
$ws.Range("K62").Value = 5975
$ws.Range("L62").Value = 17867.264
$ws.Range("M62").Value = -5351
$ws.Range("N62").Value = -19115.264

# Row 65
$ws.Range("H65").Value = 15799.044
$ws.Range("I65").Value = 5975
$ws.Range("J65").Value = 17867.264
$ws.Range("K65").Value = 29875
$ws.Range("L65").Value = 89336.31999999999
$ws.Range("M65").Value = -26755
$ws.Range("N65").Value = -95576.31999999999

# Row 113
$ws.Range("H113").Value = 1290.909
$ws.Range("I113").Value = 984.7692
$ws.Range("K113").Value = 2954.3076
$ws.Range("M113").Value = -784.3076000000001

# Row 122
$ws.Range("H122").Value = 5183.4
$ws.Range("I122").Value = 2531.625
$ws.Range("J122").Value = 8214
$ws.Range("K122").Value = 7594.875
$ws.Range("L122").Value = 24642
$ws.Range("M122").Value = -5144.875
$ws.Range("N122").Value = -29542

# Row 132
$ws.Range("H132").Value = 2214.6416
$ws.Range("I132").Value = 2007.4445
$ws.Range("K132").Value = 6022.333500000001
$ws.Range("M132").Value = -3492.333500000001

# Row 134
$ws.Range("H134").Value = 82949.5
$ws.Range("J134").Value = 82949.5
$ws.Range("L134").Value = 248848.5
$ws.Range("N134").Value = -253918.5

# Row 136
$ws.Range("H136").Value = 1485.7959
$ws.Range("I136").Value = 1020
$ws.Range("K136").Value = 3060
$ws.Range("M136").Value = -510

